# Apply the "polished, make all deliverables" update to PurchaseList.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PurchaseList")

# Update the printed date/time stamp (shared strings used in row 33)
$ws.Range("G33").Value = "19/10/2015"
$ws.Range("H33").Value = "10:42:12"

# Update Supplier Stock 1 (column L) quantities
$ws.Range("L5").Value  = 4226408
$ws.Range("L9").Value  = 197098
$ws.Range("L12").Value = 8920
$ws.Range("L13").Value = 9321
$ws.Range("L14").Value = 13567
$ws.Range("L17").Value = 1575457
$ws.Range("L18").Value = 48690
$ws.Range("L19").Value = 61778
$ws.Range("L20").Value = 1844557
$ws.Range("L22").Value = 1232071

# Refresh the "report created" NOW() timestamp cached value
$ws.Range("M33").Formula = "=NOW()"
